$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "pontos.(Nota final+P_recuperação)/2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "pontos.^l^l(Nota final+P_recuperação)/2",
    2
)

$d.Content.Find.Execute(
    "nota final.(Nota final+P_recuperação)/2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nota final.^l^l(Nota final+P_recuperação)/2",
    2
)

$d.Content.Find.Execute(
    "226p.B)MEDRONHO",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "226p.^l^lB)MEDRONHO",
    2
)
